$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated capital structure database values for rows 2 and 3 (Kuwait Oil/Gas Integrated)
$rows = @(2, 3)
foreach ($r in $rows) {
    $ws.Range("D$r").Value  = -0.25
    $ws.Range("G$r").Value  = -0.2435897435897436
    $ws.Range("H$r").Value  = -0.2435897435897436
    $ws.Range("I$r").Value  = -0.7649572649572649
    $ws.Range("J$r").Value  = -0.7649572649572649
    $ws.Range("K$r").Value  = -8.050000000000001
    $ws.Range("L$r").Value  = -0.6880341880341881

    $ws.Range("U$r").Value  = 56.6
    $ws.Range("V$r").Value  = 0.8085714285714286
    $ws.Range("W$r").Value  = -0.09710494571773221
    $ws.Range("X$r").Value  = 0.06826899002817814
    $ws.Range("Y$r").Value  = -0.1653739357459104
    $ws.Range("Z$r").Value  = 0.3660826032540675
    $ws.Range("AA$r").Value = -0.280037546933667
    $ws.Range("AB$r").Value = 0.06768507597785188
    $ws.Range("AC$r").Value = -0.3477226229115189
    $ws.Range("AD$r").Value = 1.01
    $ws.Range("AF$r").Value = 1.01
    $ws.Range("AG$r").Value = -55.59
    $ws.Range("AH$r").Value = 0.01422334882410928
    $ws.Range("AI$r").Value = 0.01110988890111099
    $ws.Range("AJ$r").Value = -3.857737682165164
    $ws.Range("AK$r").Value = -1.620227338968231
    $ws.Range("AL$r").Value = 0.137
    $ws.Range("AM$r").Value = -1.023
    $ws.Range("AN$r").Value = -0.1444921316165951
    $ws.Range("AO$r").Value = -65.32846715328466
    $ws.Range("AP$r").Value = 7.952789699570816
    $ws.Range("AQ$r").Value = 8.748778103616813
}
